$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "22.359.89"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.43%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.566.24"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.33%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.21%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.003"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.13%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "287.13"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.15%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3777"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.29%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3270"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.05%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.48"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -5.35%  "

$ws.Range("E10").Value = "  +1.20%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07419"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.21%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.001"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.20%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.43"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.01%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.861"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.13%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.809"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.68%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.551.49"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.46%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001095"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.31%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06727"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.34%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "86.01"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.41%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.004"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.22%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.368"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.56%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "16.28"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.98%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.70"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.47%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "22.365.53"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.42%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.299"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.53%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.523"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.71%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "150.76"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.95%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.42"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.17%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.901"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.19%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "123.29"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.83%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.732.94"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.11%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.045"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.58%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.927"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.96%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.912"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.28%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.489"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.09%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08257"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.22%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02378"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.17%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06302"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.72%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2183"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.46%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.272"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.19%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.267"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.60%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.04"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.07%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6084"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.10%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.001"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.06%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.68"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.34%  "

$ws.Range("E46").Value = "  -0.14%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5893"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.01%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.001"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.31%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "124.05"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.46%  "

$ws.Range("E50").Value = "  -3.44%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07136"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.29%  "
